$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Stash existing E-column fill styles (7, 12, 14, 15) into far-away helper cells
#     before touching any data, so we can re-apply exact formatting later regardless
#     of what happens to the original source cells. ---
$ws.Range("E12").Copy() | Out-Null
$ws.Range("ZZ1").PasteSpecial(-4122) | Out-Null
$ws.Range("E28").Copy() | Out-Null
$ws.Range("ZZ2").PasteSpecial(-4122) | Out-Null
$ws.Range("E3").Copy() | Out-Null
$ws.Range("ZZ3").PasteSpecial(-4122) | Out-Null
$ws.Range("E26").Copy() | Out-Null
$ws.Range("ZZ4").PasteSpecial(-4122) | Out-Null
# Also stash the column A / C / D base cell styles (wrapText / center / numFmt).
$ws.Range("A13").Copy() | Out-Null
$ws.Range("ZZ5").PasteSpecial(-4122) | Out-Null
$ws.Range("C13").Copy() | Out-Null
$ws.Range("ZZ6").PasteSpecial(-4122) | Out-Null
$ws.Range("D13").Copy() | Out-Null
$ws.Range("ZZ7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Extend the Excel Table (Tableau1) to the new range first, so the newly written
#     rows below are recognised as part of the table. ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E39")) | Out-Null

# --- Write row data for rows 13-39 (task review backlog, incl. 9 new Code Review tasks) ---
# Row 13: mettre a jour contrôleur code pour ajouter un commentaire (ajout à la vue d'un code)
$ws.Range("A13").Value = "mettre a jour contrôleur code pour ajouter un commentaire (ajout à la vue d'un code)"
$ws.Range("A13").Style = $ws.Range("A13").Style
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = 1
$ws.Range("D13").Formula = '=IFERROR(Tableau1[[#This Row],[Temps (estimation)]]/Tableau1[[#This Row],[Temps (réel)]], "")'
$ws.Range("ZZ1").Copy() | Out-Null
$ws.Range("E13").PasteSpecial(-4122) | Out-Null

# Row 14: mise à jour des commentaires
$ws.Range("A14").Value = "mise à jour des commentaires"
$ws.Range("A14").Style = $ws.Range("A13").Style
$ws.Range("B14").Value = 2
$ws.Range("C14").Value = 3
$ws.Range("D14").Formula = '=IFERROR(Tableau1[[#This Row],[Temps (estimation)]]/Tableau1[[#This Row],[Temps (réel)]], "")'
$ws.Range("ZZ3").Copy() | Out-Null
$ws.Range("E14").PasteSpecial(-4122) | Out-Null

# Row 15: suppression d'un commentaires
$ws.Range("A15").Value = "suppression d'un commentaires"
$ws.Range("A15").Style = $ws.Range("A13").Style
$ws.Range("B15").Value = 0.25
$ws.Range("C15").Value = 0.25
$ws.Range("D15").Formula = '=IFERROR(Tableau1[[#This Row],[Temps (estimation)]]/Tableau1[[#This Row],[Temps (réel)]], "")'
$ws.Range("ZZ3").Copy() | Out-Null
$ws.Range("E15").PasteSpecial(-4122) | Out-Null

# Row 16: trier l'affichage des codes
$ws.Range("A16").Value = "trier l'affichage des codes"
$ws.Range("A16").Style = $ws.Range("A13").Style
$ws.Range("B16").Value = 2
$ws.Range("C16").Value = 4
$ws.Range("D16").Formula = '=IFERROR(Tableau1[[#This Row],[Temps (estimation)]]/Tableau1[[#This Row],[Temps (réel)]], "")'
$ws.Range("ZZ1").Copy() | Out-Null
$ws.Range("E16").PasteSpecial(-4122) | Out-Null

# Row 17: confort utilisateur page login
$ws.Range("A17").Value = "confort utilisateur page login"
$ws.Range("A17").Style = $ws.Range("A13").Style
$ws.Range("B17").Value = 2
$ws.Range("C17").Value = 2
$ws.Range("D17").Formula = '=IFERROR(Tableau1[[#This Row],[Temps (estimation)]]/Tableau1[[#This Row],[Temps (réel)]], "")'
$ws.Range("ZZ1").Copy() | Out-Null
$ws.Range("E17").PasteSpecial(-4122) | Out-Null

# Row 18: confort utilisateur page codes
$ws.Range("A18").Value = "confort utilisateur page codes"
$ws.Range("A18").Style = $ws.Range("A13").Style
$ws.Range("B18").Value = 3
$ws.Range("C18").Value = 2.5
$ws.Range("D18").Formula = '=IFERROR(Tableau1[[#This Row],[Temps (estimation)]]/Tableau1[[#This Row],[Temps (réel)]], "")'
$ws.Range("ZZ3").Copy() | Out-Null
$ws.Range("E18").PasteSpecial(-4122) | Out-Null

# Row 19: confort utilisateur commentaire (saisie/modifcation/suppression
$ws.Range("A19").Value = "confort utilisateur commentaire (saisie/modifcation/suppression"
$ws.Range("A19").Style = $ws.Range("A13").Style
$ws.Range("B19").Value = 10
$ws.Range("C19").Value = 12
$ws.Range("D19").Formula = '=IFERROR(Tableau1[[#This Row],[Temps (estimation)]]/Tableau1[[#This Row],[Temps (réel)]], "")'
$ws.Range("ZZ3").Copy() | Out-Null
$ws.Range("E19").PasteSpecial(-4122) | Out-Null

# Row 20: confort utilisateur page code
$ws.Range("A20").Value = "confort utilisateur page code"
$ws.Range("A20").Style = $ws.Range("A13").Style
$ws.Range("B20").Value = 4
$ws.Range("C20").Value = 6
$ws.Range("D20").Formula = '=IFERROR(Tableau1[[#This Row],[Temps (estimation)]]/Tableau1[[#This Row],[Temps (réel)]], "")'
$ws.Range("ZZ3").Copy() | Out-Null
$ws.Range("E20").PasteSpecial(-4122) | Out-Null

# Row 21: hasher le mot de passe
$ws.Range("A21").Value = "hasher le mot de passe"
$ws.Range("A21").Style = $ws.Range("A13").Style
$ws.Range("B21").Value = 1
$ws.Range("C21").Value = 0.5
$ws.Range("D21").Formula = '=IFERROR(Tableau1[[#This Row],[Temps (estimation)]]/Tableau1[[#This Row],[Temps (réel)]], "")'
$ws.Range("ZZ1").Copy() | Out-Null
$ws.Range("E21").PasteSpecial(-4122) | Out-Null

# Row 22: créer la base de données selon le modèle définit
$ws.Range("A22").Value = "créer la base de données selon le modèle définit"
$ws.Range("A22").Style = $ws.Range("A13").Style
$ws.Range("B22").Value = 0.5
$ws.Range("C22").Value = 0.5
$ws.Range("D22").Formula = '=IFERROR(Tableau1[[#This Row],[Temps (estimation)]]/Tableau1[[#This Row],[Temps (réel)]], "")'
$ws.Range("ZZ1").Copy() | Out-Null
$ws.Range("E22").PasteSpecial(-4122) | Out-Null

# Row 23: mettre à jour le routeur
$ws.Range("A23").Value = "mettre à jour le routeur"
$ws.Range("A23").Style = $ws.Range("A13").Style
$ws.Range("B23").Value = 1
$ws.Range("C23").Value = 1
$ws.Range("D23").Formula = '=IFERROR(Tableau1[[#This Row],[Temps (estimation)]]/Tableau1[[#This Row],[Temps (réel)]], "")'
$ws.Range("ZZ3").Copy() | Out-Null
$ws.Range("E23").PasteSpecial(-4122) | Out-Null

# Row 24: adaptation pour appareil mobile
$ws.Range("A24").Value = "adaptation pour appareil mobile"
$ws.Range("A24").Style = $ws.Range("A13").Style
$ws.Range("B24").Value = 10
$ws.Range("C24").Value = 10
$ws.Range("D24").Formula = '=IFERROR(Tableau1[[#This Row],[Temps (estimation)]]/Tableau1[[#This Row],[Temps (réel)]], "")'
$ws.Range("ZZ1").Copy() | Out-Null
$ws.Range("E24").PasteSpecial(-4122) | Out-Null

# Row 25: upvote des commentaires
$ws.Range("A25").Value = "upvote des commentaires"
$ws.Range("A25").Style = $ws.Range("A13").Style
$ws.Range("B25").Value = 15
$ws.Range("C25").Value = 15
$ws.Range("D25").Formula = '=IFERROR(Tableau1[[#This Row],[Temps (estimation)]]/Tableau1[[#This Row],[Temps (réel)]], "")'
$ws.Range("ZZ4").Copy() | Out-Null
$ws.Range("E25").PasteSpecial(-4122) | Out-Null

# Row 26: menu arborescent
$ws.Range("A26").Value = "menu arborescent"
$ws.Range("A26").Style = $ws.Range("A13").Style
$ws.Range("B26").Value = 20
$ws.Range("C26").ClearContents()
$ws.Range("D26").Formula = '=IFERROR(Tableau1[[#This Row],[Temps (estimation)]]/Tableau1[[#This Row],[Temps (réel)]], "")'
$ws.Range("ZZ4").Copy() | Out-Null
$ws.Range("E26").PasteSpecial(-4122) | Out-Null

# Row 27: créer Page de SignUp
$ws.Range("A27").Value = "créer Page de SignUp"
$ws.Range("A27").Style = $ws.Range("A13").Style
$ws.Range("B27").Value = 0.5
$ws.Range("C27").Value = 0.5
$ws.Range("D27").Formula = '=IFERROR(Tableau1[[#This Row],[Temps (estimation)]]/Tableau1[[#This Row],[Temps (réel)]], "")'
$ws.Range("ZZ2").Copy() | Out-Null
$ws.Range("E27").PasteSpecial(-4122) | Out-Null

# Row 28: Ajouter des User(BDD+Controller)
$ws.Range("A28").Value = "Ajouter des User(BDD+Controller)"
$ws.Range("A28").Style = $ws.Range("A13").Style
$ws.Range("B28").Value = 2
$ws.Range("C28").Value = 2
$ws.Range("D28").Formula = '=IFERROR(Tableau1[[#This Row],[Temps (estimation)]]/Tableau1[[#This Row],[Temps (réel)]], "")'
$ws.Range("ZZ2").Copy() | Out-Null
$ws.Range("E28").PasteSpecial(-4122) | Out-Null

# Row 29: ajouter un cancel dans la page de login
$ws.Range("A29").Value = "ajouter un cancel dans la page de login"
$ws.Range("A29").Style = $ws.Range("A13").Style
$ws.Range("B29").Value = 0.5
$ws.Range("C29").Value = 0.5
$ws.Range("D29").Formula = '=IFERROR(Tableau1[[#This Row],[Temps (estimation)]]/Tableau1[[#This Row],[Temps (réel)]], "")'
$ws.Range("ZZ1").Copy() | Out-Null
$ws.Range("E29").PasteSpecial(-4122) | Out-Null

# Row 30: ajouter une mini barre de navigation
$ws.Range("A30").Value = "ajouter une mini barre de navigation"
$ws.Range("A30").Style = $ws.Range("A13").Style
$ws.Range("B30").Value = 1
$ws.Range("C30").Value = 0.75
$ws.Range("D30").Formula = '=IFERROR(Tableau1[[#This Row],[Temps (estimation)]]/Tableau1[[#This Row],[Temps (réel)]], "")'
$ws.Range("ZZ1").Copy() | Out-Null
$ws.Range("E30").PasteSpecial(-4122) | Out-Null

# Row 31: AutoDoc d'une classe
$ws.Range("A31").Value = "AutoDoc d'une classe"
$ws.Range("A31").Style = $ws.Range("A13").Style
$ws.Range("B31").Value = 0.25
$ws.Range("C31").ClearContents()
$ws.Range("D31").Formula = '=IFERROR(Tableau1[[#This Row],[Temps (estimation)]]/Tableau1[[#This Row],[Temps (réel)]], "")'
$ws.Range("ZZ1").Copy() | Out-Null
$ws.Range("E31").PasteSpecial(-4122) | Out-Null

# Row 32: TestUnitaire
$ws.Range("A32").Value = "TestUnitaire"
$ws.Range("A32").Style = $ws.Range("A13").Style
$ws.Range("B32").Value = 0.5
$ws.Range("C32").ClearContents()
$ws.Range("D32").Formula = '=IFERROR(Tableau1[[#This Row],[Temps (estimation)]]/Tableau1[[#This Row],[Temps (réel)]], "")'
$ws.Range("ZZ1").Copy() | Out-Null
$ws.Range("E32").PasteSpecial(-4122) | Out-Null

# Row 33: Manuel Utilisateur
$ws.Range("A33").Value = "Manuel Utilisateur"
$ws.Range("A33").Style = $ws.Range("A13").Style
$ws.Range("B33").Value = 2
$ws.Range("C33").ClearContents()
$ws.Range("D33").Formula = '=IFERROR(Tableau1[[#This Row],[Temps (estimation)]]/Tableau1[[#This Row],[Temps (réel)]], "")'
$ws.Range("ZZ1").Copy() | Out-Null
$ws.Range("E33").PasteSpecial(-4122) | Out-Null

# Row 34: Revue de Codes - CSS
$ws.Range("A34").Value = "Revue de Codes - CSS"
$ws.Range("A34").Style = $ws.Range("A13").Style
$ws.Range("B34").Value = 0.5
$ws.Range("C34").ClearContents()
$ws.Range("D34").Formula = '=IFERROR(Tableau1[[#This Row],[Temps (estimation)]]/Tableau1[[#This Row],[Temps (réel)]], "")'
$ws.Range("ZZ1").Copy() | Out-Null
$ws.Range("E34").PasteSpecial(-4122) | Out-Null

# Row 35: Revue de Code - JS
$ws.Range("A35").Value = "Revue de Code - JS"
$ws.Range("A35").Style = $ws.Range("A13").Style
$ws.Range("B35").Value = 0.5
$ws.Range("C35").ClearContents()
$ws.Range("D35").Formula = '=IFERROR(Tableau1[[#This Row],[Temps (estimation)]]/Tableau1[[#This Row],[Temps (réel)]], "")'
$ws.Range("ZZ1").Copy() | Out-Null
$ws.Range("E35").PasteSpecial(-4122) | Out-Null

# Row 36: Revue de Codes -HTML
$ws.Range("A36").Value = "Revue de Codes -HTML"
$ws.Range("A36").Style = $ws.Range("A13").Style
$ws.Range("B36").Value = 0.5
$ws.Range("C36").ClearContents()
$ws.Range("D36").Formula = '=IFERROR(Tableau1[[#This Row],[Temps (estimation)]]/Tableau1[[#This Row],[Temps (réel)]], "")'
$ws.Range("ZZ1").Copy() | Out-Null
$ws.Range("E36").PasteSpecial(-4122) | Out-Null

# Row 37: Revue de codes -PHP Models
$ws.Range("A37").Value = "Revue de codes -PHP Models"
$ws.Range("A37").Style = $ws.Range("A13").Style
$ws.Range("B37").Value = 1
$ws.Range("C37").ClearContents()
$ws.Range("D37").Formula = '=IFERROR(Tableau1[[#This Row],[Temps (estimation)]]/Tableau1[[#This Row],[Temps (réel)]], "")'
$ws.Range("ZZ1").Copy() | Out-Null
$ws.Range("E37").PasteSpecial(-4122) | Out-Null

# Row 38: Revues deCodes -PHP Controllers
$ws.Range("A38").Value = "Revues deCodes -PHP Controllers"
$ws.Range("A38").Style = $ws.Range("A13").Style
$ws.Range("B38").Value = 1
$ws.Range("C38").ClearContents()
$ws.Range("D38").Formula = '=IFERROR(Tableau1[[#This Row],[Temps (estimation)]]/Tableau1[[#This Row],[Temps (réel)]], "")'
$ws.Range("ZZ1").Copy() | Out-Null
$ws.Range("E38").PasteSpecial(-4122) | Out-Null

# Row 39: Revue de Codes reste
$ws.Range("A39").Value = "Revue de Codes reste"
$ws.Range("A39").Style = $ws.Range("A13").Style
$ws.Range("B39").Value = 1
$ws.Range("C39").ClearContents()
$ws.Range("D39").Formula = '=IFERROR(Tableau1[[#This Row],[Temps (estimation)]]/Tableau1[[#This Row],[Temps (réel)]], "")'
$ws.Range("ZZ1").Copy() | Out-Null
$ws.Range("E39").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# --- Re-apply the column A/C/D base cell styles to every written row (handles the newly
#     added rows 32-39, which start out with no explicit style). ---
$ws.Range("ZZ5").Copy() | Out-Null
$ws.Range("A13:A39").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ6").Copy() | Out-Null
$ws.Range("C13:C39").PasteSpecial(-4122) | Out-Null
$ws.Range("ZZ7").Copy() | Out-Null
$ws.Range("D13:D39").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 13 and row 19 wrap to a taller row height (match the source text length).
$ws.Rows("13").RowHeight = $ws.Rows("12").RowHeight
$ws.Rows("19").RowHeight = $ws.Rows("12").RowHeight

# --- Clear helper stash cells ---
$ws.Range("ZZ1:ZZ7").Clear() | Out-Null

# --- Update sheet view (scroll position / selection) to match the saved view ---
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Range("H37").Select() | Out-Null
